$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full roster of absentee names, in the final desired top-to-bottom order.
# (Row 1 keeps its existing "Absentees" header in B1.)
$names = @(
    "AANYA JAIN",
    "AARAV DUA",
    "AAYUSH GUPTA",
    "AKSHITA PURI",
    "ANEESH RAMAN",
    "ANSHIKA",
    "ARNAV SHARMA",
    "AROUSH SETH",
    "ARSHIA KHAUND",
    "ARYAN WALIA",
    "AVNI AGGARWAL",
    "BHAVYA ARORA",
    "BHAVYA SHARMA",
    "DREESHTI KAPOOR",
    "DIPIN PANDEY",
    "DEVANSH PANDEYA",
    "EKAANSH GABA",
    "ISHANI JHA",
    "IHINA ROY",
    "LAKSHAY MALHOTRA",
    "KASHIKA TAYAL",
    "JIAH BAJAJ",
    "MAHI WADHWA",
    "PANKAJ",
    "PARTH GUPTA",
    "PRATHAM SHARMA",
    "RANVEER SOLANKI",
    "RENNIE GUPTA",
    "RIHIT RAI",
    "RISHABH SINGH",
    "ROUNAK BISWAS",
    "RUDRA VIJ",
    "SANYAM MATHUR",
    "SASHVI SINGLA",
    "SHARVI SINGHAL",
    "SUMAN",
    "UNNABH BHALLA",
    "VANSHIKA ARYA",
    "YANA VIG",
    "YUVRAJ MALIK"
)

# Write serial numbers (col A) and names (col B) for every data row.
for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

# The new rows (25-41) need the same formatting (bold, thin box border,
# centered/top aligned) already used by the existing numbered rows in
# column A. Copy that formatting down in one shot.
$ws.Range("A2").Copy()
$ws.Range("A25:A41").PasteSpecial(-4122)
